$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "36.393.79"
$ws.Range("E2").Value = "  -0.11%  "

# Row 3
$ws.Range("D3").Value = "1.942.82"
$ws.Range("E3").Value = "  -2.06%  "

# Row 4
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
$ws.Range("D5").Value = "'242.43"
$ws.Range("E5").Value = "  -0.72%  "

# Row 6
$ws.Range("D6").Value = "'0.610"
$ws.Range("E6").Value = "  -2.83%  "

# Row 7
$ws.Range("E7").Value = "  -0.09%  "

# Row 8
$ws.Range("D8").Value = "'57.12"
$ws.Range("E8").Value = "  -3.23%  "

# Row 9
$ws.Range("E9").Value = "  -3.63%  "

# Row 10
$ws.Range("D10").Value = "'0.0857"
$ws.Range("E10").Value = "  +4.74%  "

# Row 11
$ws.Range("E11").Value = "  +0.38%  "

# Row 12
$ws.Range("D12").Value = "2.226.38"
$ws.Range("E12").Value = "  -2.13%  "

# Row 13
$ws.Range("D13").Value = "'0.816"
$ws.Range("E13").Value = "  -5.54%  "

# Row 14
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'21.13"
$ws.Range("E14").Value = "  -10.34%  "

# Row 15
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'13.49"
$ws.Range("E15").Value = "  -3.71%  "

# Row 16
$ws.Range("D16").Value = "'5.18"
$ws.Range("E16").Value = "  -5.30%  "

# Row 17
$ws.Range("D17").Value = "1.943.12"
$ws.Range("E17").Value = "  -2.52%  "

# Row 18
$ws.Range("D18").Value = "36.322.35"
$ws.Range("E18").Value = "  -0.31%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0875"
$ws.Range("E19").Value = "  +1.13%  "

# Row 20
$ws.Range("D20").Value = "'69.27"
$ws.Range("E20").Value = "  -1.97%  "

# Row 21
$ws.Range("D21").Value = "'228.21"
$ws.Range("E21").Value = "  -2.60%  "

# Row 22
$ws.Range("D22").Value = "'5.01"
$ws.Range("E22").Value = "  -5.76%  "

# Row 23
$ws.Range("E23").Value = "  -0.13%  "

# Row 24
$ws.Range("D24").Value = "'2.40"
$ws.Range("E24").Value = "  -7.71%  "

# Row 25
$ws.Range("D25").Value = "'2.28"
$ws.Range("E25").Value = "  -1.30%  "

# Row 26
$ws.Range("D26").Value = "'9.23"
$ws.Range("E26").Value = "  -8.95%  "

# Row 27
$ws.Range("D27").Value = "'160.88"
$ws.Range("E27").Value = "  -0.81%  "

# Row 28
$ws.Range("D28").Value = "'0.134"
$ws.Range("E28").Value = "  +2.72%  "

# Row 29
$ws.Range("E29").Value = "  -3.34%  "

# Row 30
$ws.Range("D30").Value = "'0.118"
$ws.Range("E30").Value = "  -2.10%  "

# Row 31
$ws.Range("E31").Value = "  -6.50%  "

# Row 32
$ws.Range("D32").Value = "'4.62"
$ws.Range("E32").Value = "  -6.03%  "

# Row 33
$ws.Range("D33").Value = "'0.0632"
$ws.Range("E33").Value = "  +0.57%  "

# Row 34
$ws.Range("E34").Value = "  -4.24%  "

# Row 35
$ws.Range("B35").Value = "THORChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D35").Value = "'6.14"
$ws.Range("E35").Value = "  -3.02%  "

# Row 36
$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  -0.10%  "

# Row 37
$ws.Range("E37").Value = "  +0.69%  "

# Row 38
$ws.Range("D38").Value = "'2.15"
$ws.Range("E38").Value = "  -5.34%  "

# Row 39
$ws.Range("D39").Value = "'3.06"
$ws.Range("E39").Value = "  +0.00%  "

# Row 40
$ws.Range("D40").Value = "'0.0972"
$ws.Range("E40").Value = "  +1.56%  "

# Row 41
$ws.Range("D41").Value = "'2.86"
$ws.Range("E41").Value = "  -1.59%  "

# Row 42
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.16"
$ws.Range("E42").Value = "  -6.83%  "

# Row 43
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0210"
$ws.Range("E43").Value = "  -1.85%  "

# Row 44
$ws.Range("D44").Value = "'15.65"
$ws.Range("E44").Value = "  -3.81%  "

# Row 45
$ws.Range("D45").Value = "1.343.20"
$ws.Range("E45").Value = "  -3.14%  "

# Row 46
$ws.Range("E46").Value = "  -6.51%  "

# Row 47
$ws.Range("D47").Value = "'87.13"
$ws.Range("E47").Value = "  -6.11%  "

# Row 48
$ws.Range("D48").Value = "'7.12"
$ws.Range("E48").Value = "  -5.80%  "

# Row 49
$ws.Range("E49").Value = "  -0.76%  "

# Row 50
$ws.Range("D50").Value = "'44.24"
$ws.Range("E50").Value = "  -2.25%  "

# Row 51
$ws.Range("D51").Value = "2.118.18"
$ws.Range("E51").Value = "  -2.38%  "
